# Update the "取得日時" (acquired date/time) column with the new append
# timestamp for the batch of rows that were written in this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-10 01:43:58"
$oldTimestamp = "2025-10-10 01:16:43"

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
